# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
#
# Several adjacent row-pairs in the "Saudi Arabia Division 1" sheet had their
# match data (everything except the running index in column A) swapped
# between the two rows. Re-apply the same swap: for each pair, exchange the
# contents of columns B..AC (id, teams, odds, etc.) between the two rows,
# leaving column A (the sequential row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colStart = 2   # column B
$colEnd   = 29  # column AC

$rowPairs = @(
    @(104, 105),
    @(115, 116),
    @(125, 126),
    @(181, 182),
    @(225, 226),
    @(238, 239)
)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    for ($c = $colStart; $c -le $colEnd; $c++) {
        $cellA = $ws.Cells.Item($rowA, $c)
        $cellB = $ws.Cells.Item($rowB, $c)

        $valA = $cellA.Value()
        $valB = $cellB.Value()

        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}
